$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.93
$ws.Range("B2").Value = 0.977
$ws.Range("D2").Value = 0.934
$ws.Range("E2").Value = 0.95
$ws.Range("F2").Value = 0.938
